# The source document has three paragraphs whose text was authored as one
# long run-on string ("1) ...cargas.2) Campo..." etc. with no separators
# between numbered items / references). The edit reflows each of those
# single <w:t> runs into a sequence of <w:t>...</w:t><w:br/> pairs - i.e. a
# manual line break is inserted after each numbered topic / bibliography
# entry, while keeping everything inside a single <w:r> run (so existing
# run-level formatting such as the <w:i/> italics on the English variant is
# preserved).
#
# Word represents a manual line break (<w:br/>) as the vertical-tab
# character Chr(11) inside Range.Text. Using Find.Execute to locate the
# (unique) target text and then assigning straight to Range.Text keeps
# straight apostrophes (e.g. "Coulomb's", "Gauss'") intact - going through
# Find.Execute's own Replacement text argument would smart-quote them.
$vbreak = [char]11

function Replace-WithBreaks($document, $oldText, $newText) {
    $rng = $document.Content
    $found = $rng.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find target text: $oldText"
    }
    $rng.Text = $newText
}

$d = $word.ActiveDocument

# 1) "Programa" section - Portuguese paragraph (list of 11 topics)
Replace-WithBreaks $d "1) Carga e Força elétrica: carga elétrica; condutores e isolantes; lei de Coulomb; quantização e conservação de cargas.2) Campo Elétrico: conceito; linhas de campo; carga pontual e dipolo elétrico, distribuição contínua.3) A Lei de Gauss: fluxo; aplicações em simetrias cilíndricas, planares e esféricas.4) Potencial Elétrico: conceito e cálculo; energia, potencial e campo elétrico, superfícies equipotenciais; carga puntiforme, dipolo elétrico e distribuições contínuas.5) Capacitores e Dielétricos: capacitância, energia e cálculo; associações, dielétrico.6) Corrente e Resistência Elétrica: corrente e densidade, resistência, Resistividade e Condutividade em função da temperatura; lei de Ohm, potência, semicondutores e supercondutores.7) Campos Magnéticos: lei de Biot-Savart.8) Lei de Ampère e aplicações; campo magnético de uma espira, solenoide e toroides.9) Indução Eletromagnética: conceitos; Lei de indução de Faraday; Lei de Lenz;10) Propriedades magnéticas da matéria;11) Equações de Maxwell." ("1) Carga e Força elétrica: carga elétrica; condutores e isolantes; lei de Coulomb; quantização e conservação de cargas." + $vbreak + "2) Campo Elétrico: conceito; linhas de campo; carga pontual e dipolo elétrico, distribuição contínua." + $vbreak + "3) A Lei de Gauss: fluxo; aplicações em simetrias cilíndricas, planares e esféricas." + $vbreak + "4) Potencial Elétrico: conceito e cálculo; energia, potencial e campo elétrico, superfícies equipotenciais; carga puntiforme, dipolo elétrico e distribuições contínuas." + $vbreak + "5) Capacitores e Dielétricos: capacitância, energia e cálculo; associações, dielétrico." + $vbreak + "6) Corrente e Resistência Elétrica: corrente e densidade, resistência, Resistividade e Condutividade em função da temperatura; lei de Ohm, potência, semicondutores e supercondutores." + $vbreak + "7) Campos Magnéticos: lei de Biot-Savart." + $vbreak + "8) Lei de Ampère e aplicações; campo magnético de uma espira, solenoide e toroides." + $vbreak + "9) Indução Eletromagnética: conceitos; Lei de indução de Faraday; Lei de Lenz;" + $vbreak + "10) Propriedades magnéticas da matéria;" + $vbreak + "11) Equações de Maxwell.")

# 2) "Programa" section - English (italic) paragraph (list of 11 topics)
Replace-WithBreaks $d "1) Electric charge and electric force: electric charge; conductors and insulators; Coulomb's law; quantization and conservation.2) Electric field: concepts; field lines; point charge and dipole, continuous distribution.3) Gauss' law: flow; applications in cylindrical, flat and spherical geometries.4) Electric potential: concept and calculation; energy, potential and electric field, equipotential surfaces; punctual loads, electric dipole and continuous distributions.5) Capacitors and dielectrics: capacitance, energy and calculation, associations, dielectrics.6) Electric current and resistance: current density, resistance and resistivity as a function of temperature; Ohm's law, power, semiconductors and superconductors.7) Magnetic field: Biot-Savart law.8) Ampère's law and applications: magnetic field of a coil, solenoid, and toroids.9) Electromagnetic induction and inductance: Faraday's law, Lenz's law.10) Magnetic properties of matter.11) Maxwell's equations." ("1) Electric charge and electric force: electric charge; conductors and insulators; Coulomb's law; quantization and conservation." + $vbreak + "2) Electric field: concepts; field lines; point charge and dipole, continuous distribution." + $vbreak + "3) Gauss' law: flow; applications in cylindrical, flat and spherical geometries." + $vbreak + "4) Electric potential: concept and calculation; energy, potential and electric field, equipotential surfaces; punctual loads, electric dipole and continuous distributions." + $vbreak + "5) Capacitors and dielectrics: capacitance, energy and calculation, associations, dielectrics." + $vbreak + "6) Electric current and resistance: current density, resistance and resistivity as a function of temperature; Ohm's law, power, semiconductors and superconductors." + $vbreak + "7) Magnetic field: Biot-Savart law." + $vbreak + "8) Ampère's law and applications: magnetic field of a coil, solenoid, and toroids." + $vbreak + "9) Electromagnetic induction and inductance: Faraday's law, Lenz's law." + $vbreak + "10) Magnetic properties of matter." + $vbreak + "11) Maxwell's equations.")

# 3) "Bibliografia" section paragraph (list of 5 references)
Replace-WithBreaks $d "NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 3, Edgard Blucher (2008).RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol.3, LTC (2008).TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.3, LTC (2008).SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 3, Pearson Addison Wesley (2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 3, Thomson Pioneira (2008)." ("NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 3, Edgard Blucher (2008)." + $vbreak + "RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol.3, LTC (2008)." + $vbreak + "TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.3, LTC (2008)." + $vbreak + "SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 3, Pearson Addison Wesley (2009)." + $vbreak + "JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 3, Thomson Pioneira (2008).")
